$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 336; this shifts existing rows 336:494 down to 337:495
# (Excel COM automatically carries formatting/styles down with the shifted rows)
$ws.Rows("336:336").Insert()

# Populate the newly inserted row 336 with the new weekly data record
$ws.Range("A336").Value = 9
$ws.Range("B336").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C336").Value = "Metropolitana"
$ws.Range("D336").Value = 45202
$ws.Range("E336").Value = 13
$ws.Range("F336").Value = 300000001
$ws.Range("G336").Value = "Rabanito"
$ws.Range("H336").Value = "Sin especificar"
$ws.Range("I336").Value = "Primera"
$ws.Range("J336").Value = 7000
$ws.Range("K336").Value = 3000
$ws.Range("L336").Value = 3000
$ws.Range("M336").Value = 3000
$ws.Range("N336").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O336").Value = "Provincia de Chacabuco"
$ws.Range("P336").Value = 30
$ws.Range("Q336").Value = 100
$ws.Range("R336").Value = "Hortaliza"
